# Implements: "Terminado de implementar contas para analises dos estados"
# - Change the selected month (Auxiliar!L5, linked to the Dashboard month combo box) from 6 to 7
# - Fill in the previously-empty % Cancel. (D:F) helper columns for every state row (5-31)
#   on the Auxiliar sheet, completing the per-state analysis.
# - Update the active selection on the Auxiliar sheet to match the author's final cursor spot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Auxiliar")

# Update the month selector cell that drives the whole dashboard (Julho = 7)
$ws.Range("L5").Value = 7

# Complete the per-state formulas in columns D (Total reservas), E (% Cancel.) and F (Avaliação média)
for ($r = 5; $r -le 31; $r++) {
    $ws.Range("D$r").Formula = "=COUNTIFS(Reservas[Estado],A$r,Reservas[Ano],Auxiliar!Ano,Reservas[Mês],Auxiliar!Mes)"
    $ws.Range("E$r").Formula = "=IFERROR(C$r/D$r,`"0%`")"
    $ws.Range("F$r").Formula = "=IFERROR(AVERAGEIFS(Reservas[Avaliação],Reservas[Estado],A$r,Reservas[Ano],Ano,Reservas[Mês],Auxiliar!Mes),0)"
}

# Match the author's final selection on the Auxiliar sheet
$ws.Range("H8").Select()
